$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D must keep their "text" storage type even though the new
# values often look like plain numbers (e.g. "11.75"). We force the number
# format to Text before assigning, then restore the default "Normal" style so
# the cell formatting matches the rest of the sheet.
$priceCells = @{
    "D2" = "29.424.24"
    "D3" = "1.900.76"
    "D5" = "325.25"
    "D7" = "0.4802"
    "D9" = "0.08065"
    "D11" = "23.37"
    "D12" = "1.905.97"
    "D13" = "5.947"
    "D14" = "7.062"
    "D15" = "89.81"
    "D17" = "0.06681"
    "D19" = "17.60"
    "D21" = "29.423.68"
    "D22" = "5.530"
    "D23" = "11.75"
    "D25" = "2.088.94"
    "D26" = "154.76"
    "D27" = "19.81"
    "D28" = "6.069"
    "D29" = "2.093"
    "D30" = "118.37"
    "D31" = "1.032"
    "D32" = "0.09498"
    "D33" = "1.391"
    "D34" = "3.542"
    "D36" = "0.02251"
    "D37" = "0.06068"
    "D38" = "1.177"
    "D39" = "0.5871"
    "D40" = "7.867"
    "D41" = "0.1844"
    "D42" = "10.21"
    "D43" = "2.407"
    "D44" = "1.278"
    "D45" = "0.07770"
    "D46" = "12.21"
    "D47" = "0.5521"
    "D48" = "1.920"
    "D49" = "113.47"
    "D50" = "0.2939"
    "D51" = "43.80"
}
foreach ($addr in $priceCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $priceCells[$addr]
    $rng.Style = "Normal"
}

# Remaining text cells (coin name, link, volume %) assign directly.
$textCells = @{
    "E2" = "  -0.76%  "
    "E3" = "  -0.82%  "
    "E4" = "  +0.13%  "
    "E5" = "  -2.75%  "
    "E6" = "  +0.18%  "
    "E7" = "  +2.81%  "
    "E8" = "  -1.20%  "
    "E9" = "  +0.40%  "
    "E10" = "  -1.11%  "
    "E11" = "  +4.55%  "
    "E12" = "  -3.39%  "
    "E13" = "  -0.91%  "
    "E14" = "  -1.59%  "
    "E15" = "  -0.32%  "
    "E16" = "  +0.24%  "
    "E17" = "  +1.48%  "
    "E18" = "  -0.37%  "
    "E19" = "  -1.40%  "
    "E20" = "  +0.08%  "
    "E21" = "  -0.74%  "
    "E22" = "  -0.81%  "
    "E24" = "  -2.47%  "
    "E25" = "  -5.25%  "
    "E26" = "  -0.46%  "
    "E27" = "  -0.32%  "
    "E28" = "  +5.38%  "
    "E29" = "  -2.22%  "
    "E30" = "  +0.82%  "
    "E31" = "  -3.38%  "
    "E32" = "  +0.23%  "
    "B33" = "ARBITRUM"
    "C33" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "E33" = "  -3.05%  "
    "B34" = "HuobiToken"
    "C34" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "E34" = "  -0.99%  "
    "E35" = "  -0.03%  "
    "E36" = "  -0.88%  "
    "E37" = "  -1.05%  "
    "E38" = "  -0.39%  "
    "E39" = "  -0.55%  "
    "E40" = "  -6.67%  "
    "E41" = "  -0.08%  "
    "E42" = "  -0.25%  "
    "E43" = "  +1.82%  "
    "E44" = "  +1.11%  "
    "E45" = "  +3.33%  "
    "E46" = "  +0.29%  "
    "E47" = "  -1.04%  "
    "E48" = "  -0.78%  "
    "E49" = "  +0.31%  "
    "B50" = "WOONetwork"
    "C50" = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
    "E50" = "  -1.98%  "
    "B51" = "Elrond"
    "C51" = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
    "E51" = "  -0.85%  "
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
